$wb = $excel.ActiveWorkbook

# --- Rename the two "Include ValueSets" sheets ---------------------------
$wsInclude1 = $wb.Worksheets.Item("Include ValueSets")
$wsInclude1.Name = "Include ValueSet #0"

$wsInclude2 = $wb.Worksheets.Item("Include ValueSets 2")
$wsInclude2.Name = "Include ValueSet #1"

# --- Update metadata values on the Metadata sheet -------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version: 1.1.0 -> 1.2.0
$wsMeta.Range("B3").Value = "1.2.0"

# Experimental: (blank) -> "false"
# Plain Value assignment turns a literal "false"/"true" into a real boolean
# cell, which does not match the source data (everything here is text).
# Route it through a formula -> values-only paste so it lands back as a
# genuine text (shared-string) cell instead.
$expCell = $wsMeta.Range("B7")
$expCell.Formula = '="false"'
$expCell.Copy()
$expCell.PasteSpecial(-4163)

# Date
$wsMeta.Range("B8").Value = "2024-10-31T19:21:51+01:00"

# Contact
$wsMeta.Range("B10").Value = "KL (http://www.kl.dk)"

# Jurisdiction: Denmark -> (blank)
$wsMeta.Range("B11").Value = ""
